$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 139
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H40").Value = 1603.75
$ws.Range("I40").Value = 1508.5
$ws.Range("J40").Value = 1699
$ws.Range("K40").Value = 1508.5
$ws.Range("L40").Value = 1699
$ws.Range("M40").Value = -1333.5
$ws.Range("N40").Value = -2049
$ws.Range("H107").Value = 381.2353
$ws.Range("I107").Value = 353.69232
$ws.Range("J107").Value = 470.75
$ws.Range("K107").Value = 353.69232
$ws.Range("L107").Value = 470.75
$ws.Range("M107").Value = 1566.30768
$ws.Range("N107").Value = -4310.75
$ws.Range("H113").Value = 7030.2
$ws.Range("I113").Value = 5985.143
$ws.Range("J113").Value = 9468.666999999999
$ws.Range("K113").Value = 5985.143
$ws.Range("L113").Value = 9468.666999999999
$ws.Range("M113").Value = -2731.143
$ws.Range("N113").Value = -15976.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 715.5
$ws.Range("I2").Value = 687.3333
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 687.3333
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -574.3333
$ws.Range("N2").Value = -1026
$ws.Range("H38").Value = 1425419.9
$ws.Range("I38").Value = 3313
$ws.Range("J38").Value = 2492000
$ws.Range("K38").Value = 3313
$ws.Range("L38").Value = 2492000
$ws.Range("M38").Value = -2846
$ws.Range("H95").Value = 25750
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 25750
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 25750
$ws.Range("N95").Value = -31242
$ws.Range("H116").Value = 715.5
$ws.Range("I116").Value = 687.3333
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 687.3333
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1606.6667
$ws.Range("N116").Value = -5388
$ws.Range("H132").Value = 2209.1538
$ws.Range("I132").Value = 1757.6364
$ws.Range("J132").Value = 4692.5
$ws.Range("K132").Value = 5272.9092
$ws.Range("L132").Value = 14077.5
$ws.Range("M132").Value = -2742.9092
$ws.Range("N132").Value = -19137.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 715.5
$ws.Range("I3").Value = 687.3333
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 687.3333
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -573.3333
$ws.Range("N3").Value = -1028
$ws.Range("H20").Value = 1263.6666
$ws.Range("I20").Value = 1145.5
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 1145.5
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -898.5
$ws.Range("H86").Value = 1499.6
$ws.Range("I86").Value = 1124.5
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1124.5
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1.5
$ws.Range("H89").Value = 1499.6
$ws.Range("I89").Value = 1124.5
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 5622.5
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -6.5
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 209.75
$ws.Range("I7").Value = 167.88889
$ws.Range("J7").Value = 335.33334
$ws.Range("K7").Value = 167.88889
$ws.Range("L7").Value = 335.33334
$ws.Range("M7").Value = -54.88889
$ws.Range("N7").Value = -561.33334
$ws.Range("H33").Value = 38115.117
$ws.Range("I33").Value = 8565.429
$ws.Range("J33").Value = 58799.9
$ws.Range("K33").Value = 8565.429
$ws.Range("L33").Value = 58799.9
$ws.Range("M33").Value = -8186.429
$ws.Range("N33").Value = -59557.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1878
$ws.Range("I14").Value = 1878
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5634
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -5461
$ws.Range("H51").Value = 1493
$ws.Range("I51").Value = 1493.75
$ws.Range("J51").Value = 1490
$ws.Range("K51").Value = 4481.25
$ws.Range("L51").Value = 4470
$ws.Range("M51").Value = -4021.25
$ws.Range("H80").Value = 1400
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3264
$ws.Range("H83").Value = 1400
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12600
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7920
$ws.Range("H107").Value = 749.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 749.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2248.5
$ws.Range("N107").Value = -6088.5
$ws.Range("H116").Value = 999
$ws.Range("I116").Value = 999
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2997
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 445
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3000
$ws.Range("N70").Value = -3540
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3000
$ws.Range("N73").Value = -4872
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -254
$ws.Range("H107").Value = 1301.0769
$ws.Range("I107").Value = 765.2222
$ws.Range("J107").Value = 2506.75
$ws.Range("K107").Value = 765.2222
$ws.Range("L107").Value = 2506.75
$ws.Range("M107").Value = 1154.7778
$ws.Range("N107").Value = -6346.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4624.5
$ws.Range("I7").Value = 4499
$ws.Range("J7").Value = 4750
$ws.Range("K7").Value = 4499
$ws.Range("L7").Value = 4750
$ws.Range("M7").Value = -4387
$ws.Range("N7").Value = -4974
$ws.Range("H31").Value = 9403
$ws.Range("I31").Value = 5007.5
$ws.Range("J31").Value = 12333.333
$ws.Range("K31").Value = 5007.5
$ws.Range("L31").Value = 12333.333
$ws.Range("M31").Value = -4759.5
$ws.Range("N31").Value = -12829.333
$ws.Range("H32").Value = 4008.818
$ws.Range("I32").Value = 1733
$ws.Range("J32").Value = 14250
$ws.Range("K32").Value = 1733
$ws.Range("L32").Value = 14250
$ws.Range("M32").Value = -1416
$ws.Range("N32").Value = -14884
$ws.Range("H46").Value = 2494.8572
$ws.Range("I46").Value = 828
$ws.Range("J46").Value = 3745
$ws.Range("K46").Value = 828
$ws.Range("L46").Value = 3745
$ws.Range("M46").Value = -640
$ws.Range("N46").Value = -4121
$ws.Range("H61").Value = 789.5
$ws.Range("I61").Value = 727.4
$ws.Range("J61").Value = 1100
$ws.Range("K61").Value = 727.4
$ws.Range("L61").Value = 1100
$ws.Range("M61").Value = -525.4
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H113").Value = 789.5
$ws.Range("I113").Value = 727.4
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 727.4
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1442.6
$ws.Range("H126").Value = 4624.5
$ws.Range("I126").Value = 4499
$ws.Range("J126").Value = 4750
$ws.Range("K126").Value = 13497
$ws.Range("L126").Value = 14250
$ws.Range("M126").Value = -11027
$ws.Range("N126").Value = -19190
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 952.0909
$ws.Range("I107").Value = 863.6667
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 2591.0001
$ws.Range("L107").Value = 4050
$ws.Range("M107").Value = -671.0001000000002
$ws.Range("N107").Value = -7890
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 235.09091
$ws.Range("I113").Value = 228.6
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 685.8
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 1484.2
$ws.Range("H119").Value = 40000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -49676
